$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Golang Architect / Principal Backend Architect"
$ws.Range("B32").Value = "https://www.dice.com/job-detail/4c45285a-0956-486e-bc87-dc3cdac57f5c"
$ws.Range("C32").Value = "Atlanta, Georgia"
$ws.Range("D32").Value = "Contract"
$ws.Range("E32").Value = "$80 - $85"
$ws.Range("F32").Value = "Montek System"
